# Update meter, ppc import-specification
# Splits the existing "Reactive Energy Lead/Lag" rows into explicit
# exported/imported capacitive & inductive interval-energy rows, inserting
# two new rows (shifting the STATE/ERROR/QS_TX/QS_RX rows down by two).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the existing row 41 (STATE[1..x]) so that the
# inductive-energy pair (exported/imported) and capacitive-energy pair
# (exported/imported) each get their own row, pushing everything below down
# by two rows (old 41-44 -> new 43-46).
$ws.Rows.Item(41).Resize(2).Insert()

# Row 39: capacitive exported interval energy (existing row, update description)
$ws.Range("A39").Value = "datapoints"
$ws.Range("B39").Value = "M_EQ_CAP_INT_EXP"
$ws.Range("C39").Value = "kVArh"
$ws.Range("D39").Value = "Reactive energy exported capacitively per interval"

# Row 40: capacitive imported interval energy (new row)
$ws.Range("A40").Value = "datapoints"
$ws.Range("B40").Value = "M_EQ_CAP_INT_IMP"
$ws.Range("C40").Value = "kVArh"
$ws.Range("D40").Value = "Reactive energy imported capacitively per interval"

# Row 41: inductive exported interval energy (existing row, update description)
$ws.Range("A41").Value = "datapoints"
$ws.Range("B41").Value = "M_EQ_IND_INT_EXP"
$ws.Range("C41").Value = "kVArh"
$ws.Range("D41").Value = "Reactive energy exported inductively per interval"

# Row 42: inductive imported interval energy (new row)
$ws.Range("A42").Value = "datapoints"
$ws.Range("B42").Value = "M_EQ_IND_INT_IMP"
$ws.Range("C42").Value = "kVArh"
$ws.Range("D42").Value = "Reactive energy imported inductively per interval"

# Rows 43-46 keep their previous content (STATE[1..x], ERROR[1..x], QS_TX, QS_RX)
# which the insert has shifted down from rows 41-44; just restate them for clarity.
$ws.Range("A43").Value = "datapoints"
$ws.Range("B43").Value = "STATE[1..x]"

$ws.Range("A44").Value = "datapoints"
$ws.Range("B44").Value = "ERROR[1..x]"
$ws.Range("D44").Value = "Global meter error conditions"

$ws.Range("A45").Value = "datapoints"
$ws.Range("B45").Value = "QS_TX"
$ws.Range("D45").Value = "Telegrams transmitted (communication quality)"

$ws.Range("A46").Value = "datapoints"
$ws.Range("B46").Value = "QS_RX"
$ws.Range("D46").Value = "Telegrams received (communication quality)"

$wb.Save()
